# Adding test cases for HUB left & right navigation.
#
# 1. Fix the "hight" -> "height" typo in the header row (I1) of the "hub"
#    sheet.
# 2. Add a new data row (row 7) for a "HubMenu" entry, matching the shape
#    of the existing rows (object id in column A, focus coordinates in
#    columns J/K).
# 3. Move the view/selection to the newly-added row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("hub")
$ws.Activate()

# -- fix header typo --------------------------------------------------
$ws.Range("I1").Value = "height"

# -- new row 7: HubMenu -------------------------------------------------
$ws.Range("A7").Value = "HubMenu"
$ws.Range("J7").Value = 365
$ws.Range("K7").Value = 370

# -- scroll/selection so the new row is visible --------------------------
$excel.ActiveWindow.ScrollColumn = 2
$ws.Range("J7:K7").Select()
